# Weekly fruit/vegetable price update: a new reporting date (2022-08-03,
# serial 44776) is inserted at the top of the data block. Two new rows
# ("Primera" and "Segunda" quality) are inserted at row 665, pushing all
# existing data rows (665-757) down by two (to 667-759).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 665, shifting the rest down.
$ws.Range("A665:A666").EntireRow.Insert()

# New row 665 - "Primera" quality for 2022-08-03 (serial 44776)
$ws.Range("A665").Value = 3
$ws.Range("B665").Value = 'Femacal de La Calera'
$ws.Range("C665").Value = 'Coquimbo'
$ws.Range("D665").Value = 44776
$ws.Range("E665").Value = 5
$ws.Range("F665").Value = 100112023
$ws.Range("G665").Value = 'Brócoli'
$ws.Range("H665").Value = 'Sin especificar'
$ws.Range("I665").Value = 'Primera'
$ws.Range("J665").Value = 2400
$ws.Range("K665").Value = 900
$ws.Range("L665").Value = 950
$ws.Range("M665").Value = 925
$ws.Range("N665").Value = '$/unidad'
$ws.Range("O665").Value = 'Provincia de Quillota'
$ws.Range("P665").Value = 925
$ws.Range("Q665").Value = 1
$ws.Range("R665").Value = 'Hortaliza'

# New row 666 - "Segunda" quality for 2022-08-03 (serial 44776)
$ws.Range("A666").Value = 3
$ws.Range("B666").Value = 'Femacal de La Calera'
$ws.Range("C666").Value = 'Coquimbo'
$ws.Range("D666").Value = 44776
$ws.Range("E666").Value = 5
$ws.Range("F666").Value = 100112023
$ws.Range("G666").Value = 'Brócoli'
$ws.Range("H666").Value = 'Sin especificar'
$ws.Range("I666").Value = 'Segunda'
$ws.Range("J666").Value = 1100
$ws.Range("K666").Value = 700
$ws.Range("L666").Value = 700
$ws.Range("M666").Value = 700
$ws.Range("N666").Value = '$/unidad'
$ws.Range("O666").Value = 'Provincia de Quillota'
$ws.Range("P666").Value = 700
$ws.Range("Q666").Value = 1
$ws.Range("R666").Value = 'Hortaliza'
